# "boucle Energie de ciqual_import, tests pour validation"
#
# Fill in the missing Ciqual constituent codes on the header/code row (row 1)
# for the columns that sit between the already-populated "Fibres" (U1) and
# "AG saturés" (Z1) codes:
#   V1 -> Polyols totaux           -> 34000
#   X1 -> Alcool                   -> 60000
#   Y1 -> Acides organiques        -> 65000
# (W1 / Cendres is intentionally left without a Ciqual code, as before.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V1").Value = 34000
$ws.Range("X1").Value = 60000
$ws.Range("Y1").Value = 65000

# Re-arrange the view: split/scroll the window over towards the right-hand
# side of the table (around the newly touched columns) instead of the old
# F1/L13:L14 view used while editing the left-hand columns.
$win = $excel.ActiveWindow
$win.SplitRow = 0
$win.SplitColumn = 16
$ws.Range("Q1").Select()
